$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Shift existing columns right to make room for a new "Planned Date" column:
#    old M (Notes) -> N, old L (Assigned person) -> M
$ws.Range("M1:M14").Copy($ws.Range("N1:N14"))
$ws.Range("L1:L14").Copy($ws.Range("M1:M14"))

# 2) Build the new L column ("Planned Date"). Copy K's formatting (it already
#    matches the styling the new date column should use), then overwrite values.
$ws.Range("K1:K14").Copy($ws.Range("L1:L14"))

$ws.Range("L1").Value = "Planned Date"

# 3) Fill in planned-date values for each task (column L)
$ws.Range("L2").Value = 42792
$ws.Range("L3").Value = 42805
$ws.Range("L4").Value = 42805
$ws.Range("L5").Value = 42812
$ws.Range("L6").Value = 42812
$ws.Range("L7").Value = 42819
$ws.Range("L8").Value = 42829
$ws.Range("L9").Value = 42836
$ws.Range("L10").Value = 42840
$ws.Range("L11").Value = 42849
$ws.Range("L12").Value = 42852
$ws.Range("L13").Value = 42855
$ws.Range("L14").Value = 42859

# 4) Fill in the new notes that explain the planned vs actual date (column N)
$ws.Range("N10").Value = "Pushed at the end due to optimization issued"
$ws.Range("N11").Value = "Implemented earlier "
$ws.Range("N3").Value = "Implementation of new elements"
$ws.Range("N5").Value = "Testing on others"
$ws.Range("N6").Value = "On time ~"
$ws.Range("N7").Value = "Textures took longer to adjust"
$ws.Range("N8").Value = "Had to wait for other implementations"
$ws.Range("N9").Value = "Game mechanics had to be implemented"
$ws.Range("N12").Value = "Added non-planned features (i.e. rotation)"

# 5) New column N width (closest representable value to the authored 40.7109375)
$ws.Columns("N").ColumnWidth = 39.8

# 6) Selection used by the author after the edit
$ws.Range("R10").Select()
